# Auto-generated Excel COM-interop script to apply workbook edits
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 323
$ws.Range("F3").Value = 1127
$ws.Range("F6").Value = 3384
$ws.Range("F7").Value = 61
$ws.Range("F10").Value = 770
$ws.Range("F11").Value = 590
$ws.Range("F13").Value = 147
$ws.Range("F14").Value = 651
$ws.Range("F15").Value = 1780
$ws.Range("F17").Value = 364
$ws.Range("F18").Value = 35
$ws.Range("F19").Value = 57
$ws.Range("F20").Value = 664
$ws.Range("F21").Value = 420
$ws.Range("F22").Value = 761
$ws.Range("F23").Value = 79375
$ws.Range("F24").Value = 79375
$ws.Range("F26").Value = 662
$ws.Range("F27").Value = 33710
$ws.Range("F28").Value = 33711
$ws.Range("F29").Value = 508
$ws.Range("F30").Value = 24
$ws.Range("F31").Value = 19
$ws.Range("F33").Value = 37
$ws.Range("F34").Value = 976
$ws.Range("F35").Value = 297
$ws.Range("F37").Value = 594
$ws.Range("F38").Value = 1352
$ws.Range("F39").Value = 1190
$ws.Range("F40").Value = 5469
$ws.Range("F41").Value = 770
$ws.Range("F42").Value = 451
$ws.Range("F46").Value = 388
$ws.Range("F50").Value = 46
$ws.Range("F51").Value = 6

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 24
$ws.Range("F14").Value = 1806
$ws.Range("F15").Value = 26
$ws.Range("F17").Value = 76
$ws.Range("F18").Value = 413
$ws.Range("F20").Value = 74
$ws.Range("F23").Value = 520
$ws.Range("F24").Value = 520
$ws.Range("F25").Value = 12
$ws.Range("F26").Value = 771
$ws.Range("F29").Value = 25
$ws.Range("F35").Value = 1665
$ws.Range("F36").Value = 493
$ws.Range("F39").Value = 1
$ws.Range("F40").Value = 107
$ws.Range("F41").Value = 107
$ws.Range("F44").Value = 27
$ws.Range("F47").Value = 69
$ws.Range("F49").Value = 144

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 104
$ws.Range("F4").Value = 728
$ws.Range("F5").Value = 563
$ws.Range("F6").Value = 597
$ws.Range("F7").Value = 111

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 728
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 24
$ws.Range("F8").Value = 3384
$ws.Range("F9").Value = 61
$ws.Range("F11").Value = 770
$ws.Range("F12").Value = 597
$ws.Range("F14").Value = 590
$ws.Range("F16").Value = 651
$ws.Range("F17").Value = 111
$ws.Range("F18").Value = 1780
$ws.Range("F19").Value = 26
$ws.Range("F22").Value = 76
$ws.Range("F23").Value = 36
$ws.Range("F24").Value = 57
$ws.Range("F25").Value = 664
$ws.Range("F26").Value = 413
$ws.Range("F27").Value = 420
$ws.Range("F28").Value = 74
$ws.Range("F29").Value = 79374
$ws.Range("F31").Value = 33710
$ws.Range("F32").Value = 508
$ws.Range("F33").Value = 19
$ws.Range("F35").Value = 520
$ws.Range("F36").Value = 37
$ws.Range("F37").Value = 12
$ws.Range("F38").Value = 976
$ws.Range("F41").Value = 297
$ws.Range("F43").Value = 25
$ws.Range("F44").Value = 594
$ws.Range("F45").Value = 1356
$ws.Range("F46").Value = 1190
$ws.Range("F47").Value = 770
$ws.Range("F48").Value = 1665
$ws.Range("F49").Value = 451
$ws.Range("F50").Value = 107
$ws.Range("F51").Value = 27
$ws.Range("F53").Value = 69
$ws.Range("F54").Value = 144
$ws.Range("F55").Value = 6

$ws = $wb.Worksheets.Item(2)
$ws.Range("C6").Value = "上海·Laurent Coulondre“心动巴黎”2024中国巡回音乐会"
$ws.Range("D6").Value = "汾阳路20号上海音乐学院内 上海贺绿汀音乐厅"
$ws.Range("E6").Value = "2024.04.20 19:30-04.20 21:30"
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = "不可售"
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81135"
$ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202401/wXDdS5ap1705651730828.jpeg"
